# Rename the diff-report column headers so they carry the file-format-version
# suffix ("_FV2404" / "_FV2410") instead of the generic "_old" / "_new".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldToFv2404 = @{
    "Segmentname_old"          = "Segmentname_FV2404"
    "Segmentgruppe_old"        = "Segmentgruppe_FV2404"
    "Segment_old"              = "Segment_FV2404"
    "Datenelement_old"         = "Datenelement_FV2404"
    "Segment ID_old"           = "Segment ID_FV2404"
    "Code_old"                 = "Code_FV2404"
    "Qualifier_old"            = "Qualifier_FV2404"
    "Beschreibung_old"         = "Beschreibung_FV2404"
    "Bedingungsausdruck_old"   = "Bedingungsausdruck_FV2404"
    "Bedingung_old"            = "Bedingung_FV2404"
    "Segmentname_new"          = "Segmentname_FV2410"
    "Segmentgruppe_new"        = "Segmentgruppe_FV2410"
    "Segment_new"              = "Segment_FV2410"
    "Datenelement_new"         = "Datenelement_FV2410"
    "Segment ID_new"           = "Segment ID_FV2410"
    "Code_new"                 = "Code_FV2410"
    "Qualifier_new"            = "Qualifier_FV2410"
    "Beschreibung_new"         = "Beschreibung_FV2410"
    "Bedingungsausdruck_new"   = "Bedingungsausdruck_FV2410"
    "Bedingung_new"            = "Bedingung_FV2410"
}

for ($col = 1; $col -le 21; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $current = $cell.Value2
    if ($oldToFv2404.ContainsKey($current)) {
        $cell.Value = $oldToFv2404[$current]
    }
}

# Freeze the header row.
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# Turn the used range into a proper Excel Table ("Table1"), matching the
# header row that was just relabelled.
$dataRange = $ws.Range("A1:U66")
$tbl = $ws.ListObjects.Add(1, $dataRange, 0, 1)
$tbl.Name = "Table1"
